$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header "ImportantDates" -> "ImportantDate"
$ws.Range("E1").Value = "ImportantDate"

# Unhide column C
$ws.Columns("C").Hidden = $false

# Update selection to E11
$ws.Range("E11").Select()
